# Re-generate the "haul" (Общее время) time strings with proper zero-padded
# minutes and seconds, e.g. "26 ч. 9 мин. 10 сек." -> "26 ч. 09 мин. 10 сек.".
# Hours are left unpadded; only single-digit minutes/seconds get a leading
# zero. This mirrors the upstream stats generator fix (#108).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds "Общее время" (total haul time) as text like
# "<h> ч. <m> мин. <s> сек.". Walk every used data row (row 1 is the header).
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2

    if ($val -match '^(\d+) ч\. (\d+) мин\. (\d+) сек\.$') {
        $hours = $matches[1]
        $minutes = $matches[2]
        $seconds = $matches[3]

        if ($minutes.Length -eq 1 -or $seconds.Length -eq 1) {
            $paddedMinutes = $minutes.PadLeft(2, '0')
            $paddedSeconds = $seconds.PadLeft(2, '0')
            $cell.Value = "$hours ч. $paddedMinutes мин. $paddedSeconds сек."
        }
    }
}
